# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Tue Mar 26 04:34:13 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.583.90'
$ws.Range("E2").Value = '  +5.05%  '
$ws.Range("D3").Value = '3.628.66'
$ws.Range("E3").Value = '  +4.96%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.49'
$ws.Range("E5").Value = '  +1.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '192.04'
$ws.Range("E6").Value = '  +4.27%  '
$ws.Range("E7").Value = '  +2.41%  '
$ws.Range("D8").Value = '3.622.72'
$ws.Range("E8").Value = '  +4.97%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("E10").Value = '  +2.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.668'
$ws.Range("E11").Value = '  +3.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.28'
$ws.Range("E12").Value = '  +3.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000290'
$ws.Range("E13").Value = '  +4.30%  '
$ws.Range("E14").Value = '  +5.15%  '
$ws.Range("D15").Value = '4.208.04'
$ws.Range("E15").Value = '  +5.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.74'
$ws.Range("E16").Value = '  +6.01%  '
$ws.Range("D17").Value = '3.624.43'
$ws.Range("E17").Value = '  +5.03%  '
$ws.Range("D18").Value = '70.514.39'
$ws.Range("E18").Value = '  +5.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.67'
$ws.Range("E19").Value = '  +4.89%  '
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("E21").Value = '  +4.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '487.48'
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.49'
$ws.Range("E23").Value = '  +16.62%  '
$ws.Range("E24").Value = '  -1.90%  '
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '90.98'
$ws.Range("E26").Value = '  +1.49%  '
$ws.Range("E27").Value = '  +7.00%  '
$ws.Range("E28").Value = '  +3.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.68'
$ws.Range("E29").Value = '  +6.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.06'
$ws.Range("E30").Value = '  +5.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.82'
$ws.Range("E31").Value = '  +9.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '627.48'
$ws.Range("E32").Value = '  +6.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.30'
$ws.Range("E34").Value = '  +7.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '66.29'
$ws.Range("E35").Value = '  +3.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '39.25'
$ws.Range("E36").Value = '  +7.71%  '
$ws.Range("E37").Value = '  +7.51%  '
$ws.Range("D38").Value = '0.0₃0816'
$ws.Range("E38").Value = '  +6.17%  '
$ws.Range("E39").Value = '  -1.25%  '
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("E41").Value = '  +1.47%  '
$ws.Range("D42").Value = '3.305.38'
$ws.Range("E42").Value = '  +2.97%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.15'
$ws.Range("E43").Value = '  +8.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.78'
$ws.Range("E44").Value = '  +10.19%  '
$ws.Range("E45").Value = '  +5.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.31'
$ws.Range("E46").Value = '  +3.06%  '
$ws.Range("E47").Value = '  +2.89%  '
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.16'
$ws.Range("E48").Value = '  +4.82%  '
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.75'
$ws.Range("E49").Value = '  +0.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.32'
$ws.Range("E50").Value = '  +3.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  +0.07%  '
